# The commit swaps the two theme parts of the deck: the theme that was
# serialized at ppt/theme/theme1.xml (the "Integral" / "Red Violet" theme
# used by the slide master) and the theme that was serialized at
# ppt/theme/theme2.xml (the default "Office Theme" used by the notes
# master) trade places - theme1.xml ends up holding the Office Theme
# colours and theme2.xml ends up holding the former Integral/Red Violet
# colours. The font scheme and format scheme (fills/lines/effects) are
# byte-identical between the two themes, so only the colour scheme (and,
# where settable, the theme/colour-scheme display names) actually change.
#
# PowerPoint's COM object model only exposes one editable theme for the
# whole deck - $p.SlideMaster.Theme (NotesMaster/HandoutMaster/Slide all
# resolve back to the same theme) - so we drive the swap through that
# single ThemeColorScheme, applying the colours the slide master's theme
# (theme1.xml) should end up with.

$p = $ppt.ActivePresentation

$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Best-effort: keep the display names in sync with the new palette.
$theme.Name = "Office Theme"
$colors.Name = "Office"

function RgbValue([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" colours (formerly at theme2.xml),
# now the colour scheme for the slide master's theme (theme1.xml).
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$colors.Item(1).RGB  = RgbValue 0x00 0x00 0x00   # dk1      000000
$colors.Item(2).RGB  = RgbValue 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colors.Item(3).RGB  = RgbValue 0x44 0x54 0x6A   # dk2      44546A
$colors.Item(4).RGB  = RgbValue 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colors.Item(5).RGB  = RgbValue 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colors.Item(6).RGB  = RgbValue 0xED 0x7D 0x31   # accent2  ED7D31
$colors.Item(7).RGB  = RgbValue 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colors.Item(8).RGB  = RgbValue 0xFF 0xC0 0x00   # accent4  FFC000
$colors.Item(9).RGB  = RgbValue 0x44 0x72 0xC4   # accent5  4472C4
$colors.Item(10).RGB = RgbValue 0x70 0xAD 0x47   # accent6  70AD47
$colors.Item(11).RGB = RgbValue 0x05 0x63 0xC1   # hlink    0563C1
$colors.Item(12).RGB = RgbValue 0x95 0x4F 0x72   # folHlink 954F72
